# Apply the "근간 스토리" content update:
# - Add ~20 new lines of dialogue/UI text (new shared strings)
# - Insert a new B53 label ("★ 영상 보기")
# - Add a new "★ 스포츠 체험" block (rows 79-90) describing a catch-ball
#   training scene between the doctor and the android Elim
# - Add a new dialogue block (rows 93-97) with speaker tags E/M
# - Add two more section markers ("★ ", "★ 예술 활동") further down (rows 100, 110)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label placed next to the existing "영상을 좀 볼까?" cue (row 53)
$ws.Range("B53").Value = "★ 영상 보기"

# New "스포츠 체험" (sports experience) scene
$ws.Range("B79").Value = "★ 스포츠 체험"

$ws.Range("B81").Value = "공을 던져 보라고요?"
$ws.Range("B82").Value = "아, 수류탄처럼 말입니까?"
$ws.Range("B83").Value = "아니."
$ws.Range("B84").Value = "내가 받을 수 있도록 말이야."
$ws.Range("B85").Value = "네?"
$ws.Range("B86").Value = "박사님을 향해 수류탄, 아니 공을요?"
$ws.Range("B87").Value = "그래, 캐치볼을 해보자고."
$ws.Range("B88").Value = "엘림의 손이 부들부들 떨리더니 공이 바닥으로 툭하고 떨어졌다."
$ws.Range("B89").Value = "인공 관절의 힘을 제대로 조정할 자신이 없습니다. 박사님."
$ws.Range("B90").Value = "괜찮아. 마음을 편하게 먹고 던져 봐."

# New dialogue block with speaker tags in column A
$ws.Range("A93").Value = "E"
$ws.Range("B93").Value = "간호 안드로이드로 만들어진 제겐 이 정도가 한계네요."

$ws.Range("A94").Value = "M"
$ws.Range("B94").Value = "오늘은 적당히 달리는 법을 배울 거야."

$ws.Range("A95").Value = "E"
$ws.Range("B95").Value = "?"

$ws.Range("A96").Value = "M"
$ws.Range("B96").Value = "시속 10km 정도로 달려보겠어?"

$ws.Range("B97").Value = "엘림의 "

# Further section markers
$ws.Range("B100").Value = "★ "
$ws.Range("B110").Value = "★ 예술 활동"

# Restore the selection/scroll position to where editing left off
$ws.Range("B97").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 72
$win.ScrollColumn = 1

$wb.Save()
